$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (B and C) for Title and Authors, shifting old B..T to D..V
$ws.Range("B:C").Insert()

# Remove now-unused trailing empty rows (21:25)
$ws.Range("21:25").Delete()

# Header row
$ws.Range("B1").Value2 = 'Title'
$ws.Range("C1").Value2 = 'Authors'

# Title / Authors per paper row
$ws.Range("B2").Value2 = 'CloudStrike: Chaos Engineering for Securityand Resiliency in Cloud Infrastructure'
$ws.Range("C2").Value2 = 'Kennedy A. Torkura, Muhammad I. H. Sukmana, Feng Cheng, Christoph Meinel '
$ws.Range("B3").Value2 = 'Securing cloud-based military systems with Security Chaos Engineering and Artificial Intelligence'
$ws.Range("C3").Value2 = 'Martin Bedoya, Sara Palacios, Daniel Diaz-López, Pantaleone Nespoli, Estefania Laverde, Sebastián Suárez'
$ws.Range("B4").Value2 = 'Infrastructure Security Checking Service Based on Chaos Engineering Method'
$ws.Range("C4").Value2 = 'Sabina Belyaeva, Yury Yanovich'
$ws.Range("B5").Value2 = 'Evaluating operational readiness using chaos engineering simulations on Kubernetes architecture in Big Data'
$ws.Range("C5").Value2 = 'Gautam Siwach, Adinarayana Haridas, Nagaraj Chinni'
$ws.Range("B6").Value2 = 'Enhancing DevSecOps practice with Large Language Models and Security Chaos Engineering'
$ws.Range("C6").Value2 = 'Martin Bedoya, Sara Palacios, Daniel Díaz‑López, Estefania Laverde, Pantaleone Nespoli'
$ws.Range("B7").Value2 = 'Enhancing Operational Resilience of Critical Infrastructure Processes Through Chaos Engineering'
$ws.Range("C7").Value2 = 'Panagiotis Dedousis, George Stergiopoulos, George Arampatzis, Dimitris Gritzalis'
$ws.Range("B8").Value2 = 'On the Way to Automatic Exploitation of Vulnerabilities and Validation of Systems Security through Security Chaos Engineering'
$ws.Range("C8").Value2 = 'Sara Palacios Chavarro, Pantaleone Nespoli, Daniel Díaz‑López, Yury Niño Roa'
$ws.Range("B9").Value2 = 'Measuring Resiliency of System of Systems using Chaos Engineering Experiments'
$ws.Range("C9").Value2 = 'Thomas Bailey, Patrick Marchione, Pete Swartz, Raed Salih, Michael R. Clark, Robert Denz'
$ws.Range("B10").Value2 = 'On Evaluating Self-Adaptive and Self-Healing Systems using Chaos Engineering'
$ws.Range("C10").Value2 = 'Moeen Ali Naqvi, Sehrish Malik, Merve Astekin, Leon Moonen'
$ws.Range("B11").Value2 = 'Chaos Engineering for Resilience Assessment of Digital Twins'
$ws.Range("C11").Value2 = 'Mattia Fogli, Carlo Giannelli, Filippo Poltronieri, Cesare Stefanelli, Mauro Tortonesi'
$ws.Range("B12").Value2 = 'Chaos engineering experiments in middleware systems using targeted network degradation and automatic fault injection'
$ws.Range("C12").Value2 = 'Tony Pierce, Jason Schanck, Alex Groeger, Raed Salih, Michael R. Clark'
$ws.Range("B13").Value2 = 'Automated Generation of Configurable Cloud-Native Chaos Testbeds'
$ws.Range("C13").Value2 = 'Jacopo Soldani, Antonio Brogi'
$ws.Range("B14").Value2 = 'Chaos Duck: A Tool for Automatic IoT Software Fault-Tolerance Analysis'
$ws.Range("C14").Value2 = 'Igor Zavalyshyn, Thomas Given‑Wilson, Axel Legay, Ramin Sadre, Etienne Rivière'
$ws.Range("B15").Value2 = 'CSBAuditor: Proactive Security Risk Analysis for Cloud Storage Broker Systems'
$ws.Range("C15").Value2 = 'K. A. Torkura, M. I. H. Sukmana, T. Strauss, H. Graupner, F. Cheng, C. Meinel'
$ws.Range("B16").Value2 = 'Continuous auditing and threat detection in multi-cloud infrastructure'
$ws.Range("C16").Value2 = 'K. A. Torkura, M. Sukmana, F. Cheng, C. Meinel'
$ws.Range("B17").Value2 = 'Boosting Microservice Resilience: An Evaluation of Istio’s Impact on Kubernetes Clusters Under Chaos'
$ws.Range("C17").Value2 = 'S. Singh, C. H. Muntean, S. Gupta'
$ws.Range("B18").Value2 = 'μ Chaos: Moving Chaos Engineering to IoT Devices'
$ws.Range("C18").Value2 = 'Wojciech Kalka, Tomasz Szydlo'
$ws.Range("B19").Value2 = 'Chaos Engineering: New Approaches To Security'
$ws.Range("C19").Value2 = 'Jamie Lewis, Chenxi Wang'
$ws.Range("B20").Value2 = 'Continuous Resilience Testing in AWS Environments with Advanced Fault Injection Techniques'
$ws.Range("C20").Value2 = 'Durga Praveen Devi'

# Wrap text for the one Authors cell that had it applied in the edit (row 7)
$ws.Range("C7").WrapText = $true

# Column widths for the new columns (approximate best-fit look)
$ws.Columns("A:A").ColumnWidth = 35.83
$ws.Columns("B:B").ColumnWidth = 105.5
$ws.Columns("C:C").ColumnWidth = 88.5

# View / selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C27").Select()

Write-Host "done"
